$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 33

$ws.Cells.Item($row, 1).Value = 45982
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

$ws.Cells.Item($row, 2).Value = 73
$ws.Cells.Item($row, 3).Value = 82
$ws.Cells.Item($row, 4).Value = 80
